$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header + width
$ws.Range("E1").Value = "c19_ventilated_cases"
$ws.Columns.Item(5).ColumnWidth = 22.28

# New column E values (c19_ventilated_cases) for existing rows
$ws.Range("E2").Value = 105
$ws.Range("E3").Value = 99
$ws.Range("E4").Value = 94
$ws.Range("E5").Value = 90
$ws.Range("E6").Value = 88
$ws.Range("E7").Value = 94
$ws.Range("E8").Value = 91
$ws.Range("E9").Value = 94
$ws.Range("E10").Value = 88
$ws.Range("E11").Value = 87
$ws.Range("E12").Value = 83
$ws.Range("E13").Value = 73
$ws.Range("E14").Value = 69
$ws.Range("E15").Value = 67
$ws.Range("E16").Value = 68
$ws.Range("E17").Value = 72
$ws.Range("E18").Value = 69
$ws.Range("E19").Value = 73
$ws.Range("E20").Value = 67
$ws.Range("E21").Value = 67
$ws.Range("E22").Value = 60

# Row 23 previously only had the date (A23); fill in the rest of the row
$ws.Range("B23").Value = 91
$ws.Range("C23").Value = 16
$ws.Range("D23").Value = 132
$ws.Range("E23").Value = 56

# Update view: move selection to E2 (also clears any prior topLeftCell scroll)
$ws.Range("E2").Select()
